$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 90.25
$ws.Range("I5").Value = 78.14286
$ws.Range("J5").Value = 175
$ws.Range("K5").Value = 78.14286
$ws.Range("L5").Value = 175
$ws.Range("M5").Value = 36.85714
$ws.Range("N5").Value = -405

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3550.1667
$ws.Range("I43").Value = 3133
$ws.Range("J43").Value = 3967.3333
$ws.Range("K43").Value = 3133
$ws.Range("L43").Value = 3967.3333
$ws.Range("M43").Value = -3064
$ws.Range("N43").Value = -4105.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4215.727
$ws.Range("I86").Value = 2995.6667
$ws.Range("J86").Value = 5679.8
$ws.Range("K86").Value = 2995.6667
$ws.Range("L86").Value = 5679.8
$ws.Range("M86").Value = -1872.6667
$ws.Range("N86").Value = -7925.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 4215.727
$ws.Range("I89").Value = 2995.6667
$ws.Range("J89").Value = 5679.8
$ws.Range("K89").Value = 14978.3335
$ws.Range("L89").Value = 28399
$ws.Range("M89").Value = -9362.333500000001
$ws.Range("N89").Value = -39631

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 14471.786
$ws.Range("I132").Value = 1820.6774
$ws.Range("J132").Value = 50124.91
$ws.Range("K132").Value = 5462.0322
$ws.Range("L132").Value = 150374.73
$ws.Range("M132").Value = -2932.0322
$ws.Range("N132").Value = -155434.73

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 418
$ws.Range("I4").Value = 418
$ws.Range("K4").Value = 418
$ws.Range("M4").Value = -302

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4912.2134
$ws.Range("I32").Value = 5020.3613
$ws.Range("J32").Value = 2316.6667
$ws.Range("K32").Value = 5020.3613
$ws.Range("L32").Value = 2316.6667
$ws.Range("M32").Value = -4733.3613
$ws.Range("N32").Value = -2890.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 168666.33
$ws.Range("I36").Value = 168666.33
$ws.Range("K36").Value = 168666.33
$ws.Range("M36").Value = -168320.33

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H40").Value = 30495
$ws.Range("J40").Value = 30495
$ws.Range("L40").Value = 30495
$ws.Range("N40").Value = -30847

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2246.5715
$ws.Range("I45").Value = 1245.3
$ws.Range("K45").Value = 1245.3
$ws.Range("M45").Value = -868.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 878
$ws.Range("I88").Value = 1100.4286
$ws.Range("J88").Value = 722.3
$ws.Range("K88").Value = 1100.4286
$ws.Range("L88").Value = 722.3
$ws.Range("M88").Value = -694.4286
$ws.Range("N88").Value = -1534.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 878
$ws.Range("I91").Value = 1100.4286
$ws.Range("J91").Value = 722.3
$ws.Range("K91").Value = 1100.4286
$ws.Range("L91").Value = 722.3
$ws.Range("M91").Value = 303.5714
$ws.Range("N91").Value = -3530.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1536.8334
$ws.Range("I20").Value = 1678.6364
$ws.Range("K20").Value = 1678.6364
$ws.Range("M20").Value = -1431.6364

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 55558876
$ws.Range("I86").Value = 62502612
$ws.Range("K86").Value = 62502612
$ws.Range("M86").Value = -62501489

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 55558876
$ws.Range("I89").Value = 62502612
$ws.Range("K89").Value = 312513060
$ws.Range("M89").Value = -312507444

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3525.8057
$ws.Range("I107").Value = 3032.4375
$ws.Range("J107").Value = 7472.75
$ws.Range("K107").Value = 3032.4375
$ws.Range("L107").Value = 7472.75
$ws.Range("M107").Value = -1112.4375
$ws.Range("N107").Value = -11312.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 50078
$ws.Range("J130").Value = 50078
$ws.Range("L130").Value = 50078
$ws.Range("N130").Value = -60118

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 4033.5
$ws.Range("I45").Value = 4033.5
$ws.Range("K45").Value = 4033.5
$ws.Range("M45").Value = -3440.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 6060.25
$ws.Range("J7").Value = 680.5
$ws.Range("L7").Value = 2041.5
$ws.Range("N7").Value = -2265.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2056.6667
$ws.Range("I34").Value = 1835
$ws.Range("K34").Value = 5505
$ws.Range("M34").Value = -5421

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 7166.6
$ws.Range("I56").Value = 7166.6
$ws.Range("K56").Value = 7166.6
$ws.Range("M56").Value = -6636.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 2091.6155
$ws.Range("I98").Value = 2170.7144
$ws.Range("J98").Value = 1999.3334
$ws.Range("K98").Value = 6512.1432
$ws.Range("L98").Value = 5998.0002
$ws.Range("M98").Value = -5014.1432
$ws.Range("N98").Value = -8994.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2933.4285
$ws.Range("I131").Value = 1489.4615
$ws.Range("J131").Value = 5279.875
$ws.Range("K131").Value = 4468.3845
$ws.Range("L131").Value = 15839.625
$ws.Range("M131").Value = 571.6154999999999
$ws.Range("N131").Value = -25919.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 40000
$ws.Range("J48").Value = 40000
$ws.Range("L48").Value = 40000
$ws.Range("N48").Value = -40970

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 97813.586
$ws.Range("I70").Value = 187818.33
$ws.Range("J70").Value = 7808.8335
$ws.Range("K70").Value = 187818.33
$ws.Range("L70").Value = 7808.8335
$ws.Range("M70").Value = -187548.33
$ws.Range("N70").Value = -8348.833500000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 97813.586
$ws.Range("I73").Value = 187818.33
$ws.Range("J73").Value = 7808.8335
$ws.Range("K73").Value = 187818.33
$ws.Range("L73").Value = 7808.8335
$ws.Range("M73").Value = -186882.33
$ws.Range("N73").Value = -9680.833500000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 30830.5
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 30830.5
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 30830.5
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -32182.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2433.0952
$ws.Range("I132").Value = 2178.2222
$ws.Range("K132").Value = 6534.6666
$ws.Range("M132").Value = -4004.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2388.543
$ws.Range("I136").Value = 2224.7856
$ws.Range("K136").Value = 6674.3568
$ws.Range("M136").Value = -4124.3568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H59").Value = 38333.332
$ws.Range("J59").Value = 38333.332
$ws.Range("L59").Value = 38333.332
$ws.Range("N59").Value = -39809.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 977.5454999999999
$ws.Range("I100").Value = 1268.1666
$ws.Range("K100").Value = 2536.3332
$ws.Range("M100").Value = -1995.3332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2408.9285
$ws.Range("I132").Value = 1381.7778
$ws.Range("J132").Value = 4257.8
$ws.Range("K132").Value = 4145.3334
$ws.Range("L132").Value = 12773.4
$ws.Range("M132").Value = -1615.3334
$ws.Range("N132").Value = -17833.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1305.1666
$ws.Range("I136").Value = 1123.1177
$ws.Range("K136").Value = 3369.3531
$ws.Range("M136").Value = -819.3531000000003
